# looting system 1차 완성 - add "dropMonster" column before "dropChance"
# and fill in the "type" column (Sword) for all item rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "dropChance" column (E) one cell right to F, then
# write the new "dropMonster" header/data into E. (Column widths in the
# sheet's <cols> metadata are keyed to columns 5-7 and must stay put, so
# this is done as plain value moves rather than a true column insert.)
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value2 = $ws.Range("E1").Value2
$ws.Range("F2").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("E3").Value2
$ws.Range("F4").Value2 = $ws.Range("E4").Value2
$ws.Range("F5").Value2 = $ws.Range("E5").Value2

# New header + data for the "dropMonster" column
$ws.Range("E1").Value2 = "dropMonster"
$ws.Range("E2").Value2 = "Orc"
$ws.Range("E3").Value2 = "Goblin"
$ws.Range("E4").Value2 = "Orc"
$ws.Range("E5").Value2 = "Orc"

# Fill in "type" column (C) for all item rows with "Sword"
$ws.Range("C2").Value2 = "Sword"
$ws.Range("C3").Value2 = "Sword"
$ws.Range("C4").Value2 = "Sword"
$ws.Range("C5").Value2 = "Sword"

# Update the active selection as in the saved workbook
$ws.Range("C7").Select()
